$d = $word.ActiveDocument

$replacements = @(
    @{old="113×7="; new="469×4="},
    @{old="139×7="; new="383×2="},
    @{old="734×9="; new="754×5="},
    @{old="503×7="; new="902×8="},
    @{old="550×4="; new="817×9="},
    @{old="839×7="; new="644×3="},
    @{old="605×4="; new="161×9="},
    @{old="266×5="; new="591×3="},
    @{old="839×9="; new="541×8="},
    @{old="887×5="; new="975×4="},
    @{old="104×4="; new="133×5="},
    @{old="165×2="; new="251×9="},
    @{old="583×3="; new="915×7="},
    @{old="476×3="; new="491×8="},
    @{old="990×9="; new="193×9="},
    @{old="587×2="; new="817×4="},
    @{old="439×8="; new="893×9="},
    @{old="399×7="; new="810×7="},
    @{old="114×6="; new="185×3="},
    @{old="529×2="; new="160×7="},
    @{old="632×8="; new="617×8="},
    @{old="736×3="; new="880×4="},
    @{old="518×8="; new="482×4="},
    @{old="842×8="; new="633×9="},
    @{old="136×9="; new="423×8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
